$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Default green"
$ws.Range("C1").Value = "Green"
$ws.Range("D1").Value = "Yellow"
$ws.Range("E1").Value = "Orange"
$ws.Range("F1").Value = "Brown"
$ws.Range("G1").Value = "Red"
$ws.Range("H1").Value = "Default Red"
$ws.Range("I1").Value = "Blue"
$ws.Range("J1").Value = ""
